$wb = $excel.ActiveWorkbook

$updates = @{
    "F2" = 632
    "F3" = 478
    "F8" = 1323
    "F9" = 3964
    "F10" = 86
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
